{"js": "// 1) Contact line: \" \" + \"\u2022 +91 7767880235 \u2022 \" (two adjacent runs with\n//    identical formatting) collapse into a single run \" \u2022 +91 7767880235 \u2022 \".\n//    Re-inserting the same visible text over the matched range merges the\n//    runs it spans into one, using the first run's formatting - exactly\n//    the structural change the diff shows (no visible text change).\nconst contactSearch = context.document.body.search(\" \u2022 +91 7767880235 \u2022 \", { matchCase: true });\ncontactSearch.load(\"items\");\nawait context.sync();\nif (contactSearch.items.length > 0) {\n  contactSearch.items[0].insertText(\" \u2022 +91 7767880235 \u2022 \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Summary paragraph: replace its text while keeping the paragraph's\n//    run formatting (rFonts/sz) intact.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst oldSummary = \"Passionate GenAI Engineer with 4+ years of experience in developing AI solutions using Python, GPT-4, and LLMs, specializing in health tech and regulatory intelligence. Proven expertise in deploying models with Docker and Kubernetes, and creating pipelines with Airflow. Granted patent holder, published author, and PharmaSUG presenter. Actively exploring the intersection of AI and life sciences through academic and industry-led innovation. \";\nconst newSummary = \"Passionate GenAI Engineer with 4+ years of experience managing full lifecycle of AI/ML and GenAI projects, collaborating with cross-functional teams, and developing AI solutions in health tech and regulatory intelligence. Granted patent holder, published author, and PharmaSUG presenter. Collaborated with cross-functional teams to explore the intersection of AI and life sciences, ensuring quality, scalability, and business alignment of AI solutions. \";\n\nconst oldSkills = \"Python, Data Science, Machine Learning, NLP, Generative AI, GPT-4, Lang Chain, Fine tuning LLMs, Docker, Kubernetes, Airflow, TensorFlow, PyTorch, Keras, RAG, Azure, OpenAI, REST APIs, AWS (Bedrock), RASA \";\nconst newSkills = \"Python, Data Science, Machine Learning, NLP, Generative AI, RAG, Azure, OpenAI, REST APIs, AWS (Bedrock), RASA, Docker, Kubernetes, Airflow, TensorFlow, PyTorch, Keras \";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text === oldSummary) {\n    paragraphs.items[i].getRange().insertText(newSummary, Word.InsertLocation.replace);\n  } else if (text === oldSkills) {\n    paragraphs.items[i].getRange().insertText(newSkills, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Contact line: the lone space run \" \" immediately followed by the\n#    \"\u2022 +91 7767880235 \u2022 \" run (identical formatting) collapse into a\n#    single run \" \u2022 +91 7767880235 \u2022 \". Locate the phone number, then\n#    widen the range by the fixed, literal text immediately surrounding\n#    it (\" \u2022 +91 \" before / \" \u2022 \" after) to land exactly on the two\n#    runs that must merge.\n$bullet = [char]8226\n$phoneFind = $d.Content\n$find = $phoneFind.Find\n$find.ClearFormatting()\n$find.Text = \"7767880235\"\n$found = $find.Execute()\nif ($found) {\n  $prefix = \" \" + $bullet + \" +91 \"\n  $suffix = \" \" + $bullet + \" \"\n  $mergeStart = $phoneFind.Start - $prefix.Length\n  $mergeEnd = $phoneFind.End + $suffix.Length\n  $mergeRange = $d.Range($mergeStart, $mergeEnd)\n  $finalText = $prefix + \"7767880235\" + $suffix\n  # Setting identical text is treated as a no-op by the engine (it will\n  # not actually restructure/merge the underlying runs), so first punch\n  # in a placeholder to force a real content mutation, then write the\n  # real text back over the freshly-written (single-run) range.\n  $placeholder = \"@@MERGE_PLACEHOLDER@@\"\n  $mergeRange.Text = $placeholder\n  $mergeRange2 = $d.Range($mergeStart, $mergeStart + $placeholder.Length)\n  $mergeRange2.Text = $finalText\n}\n\n# 2) Summary paragraph: replace the whole sentence, preserving the\n#    paragraph's existing run formatting.\n$oldSummary = \"Passionate GenAI Engineer with 4+ years of experience in developing AI solutions using Python, GPT-4, and LLMs, specializing in health tech and regulatory intelligence. Proven expertise in deploying models with Docker and Kubernetes, and creating pipelines with Airflow. Granted patent holder, published author, and PharmaSUG presenter. Actively exploring the intersection of AI and life sciences through academic and industry-led innovation. \"\n$newSummary = \"Passionate GenAI Engineer with 4+ years of experience managing full lifecycle of AI/ML and GenAI projects, collaborating with cross-functional teams, and developing AI solutions in health tech and regulatory intelligence. Granted patent holder, published author, and PharmaSUG presenter. Collaborated with cross-functional teams to explore the intersection of AI and life sciences, ensuring quality, scalability, and business alignment of AI solutions. \"\n\n$summaryRng = $d.Content\n$summaryFind = $summaryRng.Find\n$summaryFind.ClearFormatting()\n$summaryFind.Text = $oldSummary\nif ($summaryFind.Execute()) {\n  $summaryRng.Text = $newSummary\n}\n\n# 3) Technical skills line: replace the whole sentence (skills\n#    reordered / trimmed), preserving the paragraph's run formatting.\n$oldSkills = \"Python, Data Science, Machine Learning, NLP, Generative AI, GPT-4, Lang Chain, Fine tuning LLMs, Docker, Kubernetes, Airflow, TensorFlow, PyTorch, Keras, RAG, Azure, OpenAI, REST APIs, AWS (Bedrock), RASA \"\n$newSkills = \"Python, Data Science, Machine Learning, NLP, Generative AI, RAG, Azure, OpenAI, REST APIs, AWS (Bedrock), RASA, Docker, Kubernetes, Airflow, TensorFlow, PyTorch, Keras \"\n\n$skillsRng = $d.Content\n$skillsFind = $skillsRng.Find\n$skillsFind.ClearFormatting()\n$skillsFind.Text = $oldSkills\nif ($skillsFind.Execute()) {\n  $skillsRng.Text = $newSkills\n}\n"}
